$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a "last changed" date for every data row
# (rows 2-351). Bump the serial date value from 45192 (2023-09-23) to
# 45202 (2023-10-03) for every one of those rows.
$ws.Range("C2:C351").Value = 45202
